$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: fill in the Screenshot Sent status in column I
$ws.Range("I3").Value = "Screenshot Sent"

# Row 4: new ticket entry
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "ClassPlus"
$ws.Range("C4").Value = 1285451
$ws.Range("D3").Copy($ws.Range("D4"))
$ws.Range("E4").Value = "Anirban Chakraborty"
$ws.Range("F4").Value = "Cannot Login to App"
$ws.Range("G4").Value = "Pending"
$ws.Range("I4").Value = "Video Sent"

# Column I now holds data and needs a best-fit width like the other data columns
$ws.Columns.Item(9).ColumnWidth = 15

# Update the active selection to I4
$ws.Range("I4").Select()
